# Apply the "Gen -> MaxFES" rework:
#  - Rename header A1 and rescale column A from generation counts to MaxFES fractions
#  - Replace the "Run 50" column (AZ) with a true Mean-of-50-runs column
#  - Drop the old trailing "Mean" column (BA), which shifts the used range back to AZ

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header + column A (MaxFES axis) ---------------------------------------
$ws.Range("A1").Value = "MaxFES"

$ws.Range("A3").Value  = 0.001
$ws.Range("A4").Value  = 0.01
$ws.Range("A5").Value  = 0.1
$ws.Range("A6").Value  = 0.2
$ws.Range("A7").Value  = 0.3
$ws.Range("A8").Value  = 0.4
$ws.Range("A9").Value  = 0.5
$ws.Range("A10").Value = 0.6
$ws.Range("A11").Value = 0.7
$ws.Range("A12").Value = 0.8
$ws.Range("A13").Value = 0.9
$ws.Range("A14").Value = 1

# --- AZ column becomes the Mean of Run 0..Run 49 (columns B:AY) ------------
$ws.Range("AZ1").Value  = "Mean"
$ws.Range("AZ2").Value  = 62.98233691
$ws.Range("AZ3").Value  = 48.35530898
$ws.Range("AZ4").Value  = 5.88982435
$ws.Range("AZ5").Value  = 0.27429377
$ws.Range("AZ6").Value  = 0.22996758
$ws.Range("AZ7").Value  = 0.21019437
$ws.Range("AZ8").Value  = 0.19405817
$ws.Range("AZ9").Value  = 0.18077652
$ws.Range("AZ10").Value = 0.17462243
$ws.Range("AZ11").Value = 0.16740299
$ws.Range("AZ12").Value = 0.15974181
$ws.Range("AZ13").Value = 0.15348506
$ws.Range("AZ14").Value = 0.14979757

# --- Drop the old "Mean" column entirely (was BA) ---------------------------
$ws.Range("BA1:BA14").EntireColumn.Delete()
